$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.825.71"
$ws.Range("E2").Value = "  -1.12%  "

$ws.Range("D3").Value = "3.480.34"
$ws.Range("E3").Value = "  -1.09%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.88"
$ws.Range("E5").Value = "  -0.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.53"
$ws.Range("E6").Value = "  -3.62%  "

$ws.Range("D7").Value = "3.479.28"
$ws.Range("E7").Value = "  -1.12%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  -1.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.141"
$ws.Range("E10").Value = "  -2.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.90"
$ws.Range("E11").Value = "  +4.36%  "

$ws.Range("E12").Value = "  -2.92%  "

$ws.Range("D13").Value = "4.080.22"
$ws.Range("E13").Value = "  -0.80%  "

$ws.Range("E14").Value = "  -2.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.90"
$ws.Range("E15").Value = "  -3.81%  "

$ws.Range("D16").Value = "3.471.64"
$ws.Range("E16").Value = "  -1.37%  "

$ws.Range("D17").Value = "66.903.68"
$ws.Range("E17").Value = "  -1.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.62"
$ws.Range("E19").Value = "  +6.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.27"
$ws.Range("E20").Value = "  -3.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.24"
$ws.Range("E21").Value = "  -2.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "430.33"
$ws.Range("E22").Value = "  -4.72%  "

$ws.Range("E23").Value = "  -4.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.53"
$ws.Range("E24").Value = "  +0.76%  "

$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").Value = "3.616.13"
$ws.Range("E26").Value = "  -0.84%  "

$ws.Range("E27").Value = "  -6.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.68"
$ws.Range("E28").Value = "  -2.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.05"
$ws.Range("E29").Value = "  -6.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.50"
$ws.Range("E30").Value = "  -0.48%  "

$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.54"
$ws.Range("E32").Value = "  -7.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.165"
$ws.Range("E33").Value = "  -2.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.31"
$ws.Range("E34").Value = "  -1.62%  "

$ws.Range("E35").Value = "  -3.83%  "

$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.89"
$ws.Range("E37").Value = "  -1.53%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.68"
$ws.Range("E38").Value = "  -8.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "173.50"
$ws.Range("E40").Value = "  -1.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0889"
$ws.Range("E41").Value = "  -1.13%  "

$ws.Range("E42").Value = "  -2.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.890"
$ws.Range("E43").Value = "  -1.11%  "

$ws.Range("E44").Value = "  -13.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.29"
$ws.Range("E45").Value = "  -1.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "27.55"
$ws.Range("E46").Value = "  -10.60%  "

$ws.Range("E47").Value = "  -6.68%  "

$ws.Range("E48").Value = "  -4.33%  "

$ws.Range("B49").Value = "SuiNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.975"
$ws.Range("E49").Value = "  -2.50%  "

$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.36"
$ws.Range("E50").Value = "  -4.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.242"
